$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds plain-text values in the source workbook
# (e.g. "35.602.94", "1.00", "9.20"). Assigning a bare numeric-looking
# string to a General-formatted cell makes Excel coerce it to a number
# (dropping formatting like trailing zeros or thousands separators), so
# every Price write is apostrophe-prefixed to keep it stored as text --
# matching the original inlineStr cells.

$ws.Range("D2").Value = "'35.565.11"
$ws.Range("E2").Value = "  -2.20%  "

$ws.Range("D3").Value = "'1.987.75"
$ws.Range("E3").Value = "  -2.65%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'244.17"
$ws.Range("E5").Value = "  +1.84%  "

$ws.Range("D6").Value = "'0.635"
$ws.Range("E6").Value = "  -4.56%  "

$ws.Range("D7").Value = "'57.78"
$ws.Range("E7").Value = "  +8.32%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "'59.15"
$ws.Range("E9").Value = "  +1.55%  "

$ws.Range("D10").Value = "'0.361"
$ws.Range("E10").Value = "  +2.22%  "

$ws.Range("D11").Value = "'0.0732"
$ws.Range("E11").Value = "  -0.82%  "

$ws.Range("E12").Value = "  -3.30%  "

$ws.Range("D13").Value = "'0.934"
$ws.Range("E13").Value = "  +7.19%  "

$ws.Range("D14").Value = "'14.25"
$ws.Range("E14").Value = "  +0.25%  "

$ws.Range("D15").Value = "'2.277.79"
$ws.Range("E15").Value = "  -2.62%  "

$ws.Range("D16").Value = "'5.25"
$ws.Range("E16").Value = "  -0.57%  "

$ws.Range("D17").Value = "'1.995.36"
$ws.Range("E17").Value = "  -2.05%  "

$ws.Range("D18").Value = "'17.25"
$ws.Range("E18").Value = "  +8.01%  "

$ws.Range("D19").Value = "'35.591.60"
$ws.Range("E19").Value = "  -1.88%  "

$ws.Range("D20").Value = "'70.96"
$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("D21").Value = "'0.0₃0841"
$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("D22").Value = "'233.79"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").Value = "'5.15"
$ws.Range("E23").Value = "  -0.47%  "

$ws.Range("E24").Value = "  -0.13%  "

$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'2.47"
$ws.Range("E25").Value = "  +18.05%  "

$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'2.31"
$ws.Range("E26").Value = "  -1.64%  "

$ws.Range("D27").Value = "'9.20"
$ws.Range("E27").Value = "  +1.18%  "

$ws.Range("D28").Value = "'163.66"
$ws.Range("E28").Value = "  +1.33%  "

$ws.Range("D29").Value = "'19.42"
$ws.Range("E29").Value = "  -2.71%  "

$ws.Range("E30").Value = "  -1.77%  "

$ws.Range("D31").Value = "'1.14"
$ws.Range("E31").Value = "  +3.50%  "

$ws.Range("D32").Value = "'4.82"
$ws.Range("E32").Value = "  -2.11%  "

$ws.Range("D33").Value = "'0.0595"
$ws.Range("E33").Value = "  +2.53%  "

$ws.Range("D34").Value = "'0.0899"
$ws.Range("E34").Value = "  +12.34%  "

$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'2.39"
$ws.Range("E35").Value = "  +11.23%  "

$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "'4.25"
$ws.Range("E36").Value = "  -4.97%  "

$ws.Range("E37").Value = "  +0.18%  "

$ws.Range("D38").Value = "'1.79"
$ws.Range("E38").Value = "  -1.99%  "

$ws.Range("E39").Value = "  +4.02%  "

$ws.Range("D40").Value = "'1.19"
$ws.Range("E40").Value = "  -1.99%  "

$ws.Range("D41").Value = "'2.84"
$ws.Range("E41").Value = "  -0.60%  "

$ws.Range("D42").Value = "'0.0212"
$ws.Range("E42").Value = "  +0.41%  "

$ws.Range("D43").Value = "'1.09"
$ws.Range("E43").Value = "  -0.32%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "'1.389.58"
$ws.Range("E44").Value = "  +0.87%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'91.26"
$ws.Range("E45").Value = "  -0.95%  "

$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "'0.0882"
$ws.Range("E46").Value = "  -1.87%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'15.83"
$ws.Range("E47").Value = "  +3.42%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'7.49"
$ws.Range("E48").Value = "  +4.80%  "

$ws.Range("D49").Value = "'2.89"
$ws.Range("E49").Value = "  +2.37%  "

$ws.Range("D50").Value = "'2.28"
$ws.Range("E50").Value = "  +1.86%  "

$ws.Range("D51").Value = "'45.61"
$ws.Range("E51").Value = "  +6.00%  "
